$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("BP1").Value = "average_doctor_old"
$ws.Range("BQ1").Value = "average_doctor"
$ws.Range("E4").Value = 0.426
$ws.Range("F4").Value = 0.073
$ws.Range("G4").Value = 0.27
$ws.Range("N4").Value = 0.435
$ws.Range("O4").Value = 0.061
$ws.Range("P4").Value = 0.248
$ws.Range("Q4").Value = 0.025
$ws.Range("R4").Value = 0.017
$ws.Range("S4").Value = 0.132
$ws.Range("W4").Value = 0.299
$ws.Range("X4").Value = 0.11
$ws.Range("Y4").Value = 0.332
$ws.Range("AI4").Value = 0.299
$ws.Range("AJ4").Value = 0.089
$ws.Range("AK4").Value = 0.298
$ws.Range("AU4").Value = 0.191
$ws.Range("AV4").Value = 0.029
$ws.Range("AW4").Value = 0.17
$ws.Range("BA4").Value = 1.999
$ws.Range("BB4").Value = 0.157
$ws.Range("BC4").Value = 0.396
$ws.Range("BG4").Value = 0.728
$ws.Range("BH4").Value = 0.139
$ws.Range("BI4").Value = 0.373
$ws.Range("BM4").Value = 0.718
$ws.Range("BN4").Value = 0.076
$ws.Range("BO4").Value = 0.275
$ws.Range("BP4").Value = 0.666
$ws.Range("BQ4").Value = 0.708
$ws.Range("E5").Value = 0.551
$ws.Range("F5").Value = 0.08500000000000001
$ws.Range("G5").Value = 0.291
$ws.Range("N5").Value = 0.752
$ws.Range("O5").Value = 0.075
$ws.Range("P5").Value = 0.274
$ws.Range("Q5").Value = 0.016
$ws.Range("R5").Value = 0.007
$ws.Range("S5").Value = 0.08400000000000001
$ws.Range("W5").Value = 0.288
$ws.Range("X5").Value = 0.11
$ws.Range("Y5").Value = 0.332
$ws.Range("AI5").Value = 0.323
$ws.Range("AJ5").Value = 0.1
$ws.Range("AK5").Value = 0.315
$ws.Range("AU5").Value = 0.375
$ws.Range("AV5").Value = 0.095
$ws.Range("AW5").Value = 0.309
$ws.Range("BA5").Value = 1.357
$ws.Range("BB5").Value = 0.081
$ws.Range("BG5").Value = 0.401
$ws.Range("BH5").Value = 0.051
$ws.Range("BI5").Value = 0.225
$ws.Range("BM5").Value = 0.5600000000000001
$ws.Range("BN5").Value = 0.062
$ws.Range("BO5").Value = 0.249
$ws.Range("BP5").Value = 0.452
$ws.Range("BQ5").Value = 0.461
$ws.Range("E6").Value = 0.481
$ws.Range("N6").Value = 0.551
$ws.Range("Q6").Value = 0.02
$ws.Range("W6").Value = 0.293
$ws.Range("AI6").Value = 0.311
$ws.Range("AU6").Value = 0.253
$ws.Range("BA6").Value = 1.608
$ws.Range("BG6").Value = 0.517
$ws.Range("BM6").Value = 0.629
$ws.Range("BP6").Value = 0.536
$ws.Range("BQ6").Value = 0.555
$ws.Range("E7").Value = 0.52
$ws.Range("N7").Value = 0.656
$ws.Range("Q7").Value = 0.017
$ws.Range("W7").Value = 0.29
$ws.Range("AI7").Value = 0.318
$ws.Range("AU7").Value = 0.314
$ws.Range("BA7").Value = 1.447
$ws.Range("BG7").Value = 0.441
$ws.Range("BM7").Value = 0.586
$ws.Range("BP7").Value = 0.482
$ws.Range("BQ7").Value = 0.495
$ws.Range("E8").Value = 0.605
$ws.Range("F8").Value = 0.113
$ws.Range("G8").Value = 0.336
$ws.Range("N8").Value = 0.787
$ws.Range("O8").Value = 0.059
$ws.Range("P8").Value = 0.242
$ws.Range("Q8").Value = 0.018
$ws.Range("S8").Value = 0.11
$ws.Range("W8").Value = 0.318
$ws.Range("Y8").Value = 0.349
$ws.Range("AI8").Value = 0.347
$ws.Range("AJ8").Value = 0.131
$ws.Range("AK8").Value = 0.361
$ws.Range("AU8").Value = 0.319
$ws.Range("AW8").Value = 0.294
$ws.Range("BA8").Value = 1.759
$ws.Range("BB8").Value = 0.125
$ws.Range("BC8").Value = 0.354
$ws.Range("BG8").Value = 0.57
$ws.Range("BH8").Value = 0.106
$ws.Range("BI8").Value = 0.326
$ws.Range("BM8").Value = 0.7
$ws.Range("BN8").Value = 0.063
$ws.Range("BO8").Value = 0.251
$ws.Range("BP8").Value = 0.586
$ws.Range("BQ8").Value = 0.606
$ws.Range("E9").Value = 0.539
$ws.Range("N9").Value = 0.6850000000000001
$ws.Range("O9").Value = 0.216
$ws.Range("P9").Value = 0.464
$ws.Range("W9").Value = 0.213
$ws.Range("X9").Value = 0.168
$ws.Range("Y9").Value = 0.41
$ws.Range("AI9").Value = 0.27
$ws.Range("AJ9").Value = 0.197
$ws.Range("AK9").Value = 0.444
$ws.Range("BA9").Value = 1.719
$ws.Range("BB9").Value = 0.247
$ws.Range("BC9").Value = 0.497
$ws.Range("BG9").Value = 0.607
$ws.Range("BH9").Value = 0.239
$ws.Range("BI9").Value = 0.488
$ws.Range("BM9").Value = 0.663
$ws.Range("BN9").Value = 0.223
$ws.Range("BO9").Value = 0.473
$ws.Range("BP9").Value = 0.573
$ws.Range("BQ9").Value = 0.586
$ws.Range("E10").Value = 0.674
$ws.Range("F10").Value = 0.22
$ws.Range("G10").Value = 0.469
$ws.Range("N10").Value = 0.888
$ws.Range("O10").Value = 0.1
$ws.Range("P10").Value = 0.316
$ws.Range("W10").Value = 0.393
$ws.Range("X10").Value = 0.239
$ws.Range("Y10").Value = 0.488
$ws.Range("AI10").Value = 0.382
$ws.Range("AJ10").Value = 0.236
$ws.Range("AK10").Value = 0.486
$ws.Range("AU10").Value = 0.315
$ws.Range("AV10").Value = 0.216
$ws.Range("AW10").Value = 0.464
$ws.Range("BA10").Value = 2.101
$ws.Range("BB10").Value = 0.243
$ws.Range("BC10").Value = 0.493
$ws.Range("BG10").Value = 0.663
$ws.Range("BH10").Value = 0.223
$ws.Range("BI10").Value = 0.473
$ws.Range("BM10").Value = 0.854
$ws.Range("BN10").Value = 0.125
$ws.Range("BO10").Value = 0.353
$ws.Range("BP10").Value = 0.7
$ws.Range("BQ10").Value = 0.726
$ws.Range("E11").Value = 0.708
$ws.Range("F11").Value = 0.207
$ws.Range("G11").Value = 0.455
$ws.Range("N11").Value = 0.91
$ws.Range("O11").Value = 0.082
$ws.Range("P11").Value = 0.286
$ws.Range("W11").Value = 0.393
$ws.Range("X11").Value = 0.239
$ws.Range("Y11").Value = 0.488
$ws.Range("AI11").Value = 0.416
$ws.Range("AJ11").Value = 0.243
$ws.Range("AK11").Value = 0.493
$ws.Range("AU11").Value = 0.449
$ws.Range("AV11").Value = 0.247
$ws.Range("AW11").Value = 0.497
$ws.Range("BA11").Value = 2.101
$ws.Range("BB11").Value = 0.243
$ws.Range("BC11").Value = 0.493
$ws.Range("BG11").Value = 0.663
$ws.Range("BH11").Value = 0.223
$ws.Range("BI11").Value = 0.473
$ws.Range("BM11").Value = 0.854
$ws.Range("BN11").Value = 0.125
$ws.Range("BO11").Value = 0.353
$ws.Range("BP11").Value = 0.7
$ws.Range("BQ11").Value = 0.729
$ws.Range("E12").Value = 1.429
$ws.Range("F12").Value = 0.785
$ws.Range("G12").Value = 0.886
$ws.Range("N12").Value = 1.482
$ws.Range("O12").Value = 1.069
$ws.Range("P12").Value = 1.034
$ws.Range("W12").Value = 1.629
$ws.Range("X12").Value = 0.576
$ws.Range("Y12").Value = 0.759
$ws.Range("AI12").Value = 1.703
$ws.Range("AJ12").Value = 1.29
$ws.Range("AK12").Value = 1.136
$ws.Range("AU12").Value = 2.714
$ws.Range("AV12").Value = 2.68
$ws.Range("AW12").Value = 1.637
$ws.Range("BA12").Value = 3.718
$ws.Range("BB12").Value = 0.412
$ws.Range("BC12").Value = 0.642
$ws.Range("BG12").Value = 1.102
$ws.Range("BH12").Value = 0.125
$ws.Range("BI12").Value = 0.354
$ws.Range("BM12").Value = 1.289
$ws.Range("BN12").Value = 0.337
$ws.Range("BO12").Value = 0.581
$ws.Range("BP12").Value = 1.239
$ws.Range("BQ12").Value = 1.261
$ws.Range("E13").Value = 1.611
$ws.Range("F13").Value = 0.659
$ws.Range("G13").Value = 0.8120000000000001
$ws.Range("N13").Value = 2.086
$ws.Range("O13").Value = 0.9409999999999999
$ws.Range("P13").Value = 0.97
$ws.Range("W13").Value = 1.05
$ws.Range("X13").Value = 0.19
$ws.Range("Y13").Value = 0.436
$ws.Range("AI13").Value = 1.288
$ws.Range("AJ13").Value = 0.37
$ws.Range("AK13").Value = 0.608
$ws.Range("AU13").Value = 2.315
$ws.Range("AV13").Value = 0.929
$ws.Range("AW13").Value = 0.964
$ws.Range("BA13").Value = 2.386
$ws.Range("BB13").Value = 0.302
$ws.Range("BC13").Value = 0.55
$ws.Range("BG13").Value = 0.593
$ws.Range("BH13").Value = 0.07199999999999999
$ws.Range("BI13").Value = 0.268
$ws.Range("BM13").Value = 0.907
$ws.Range("BN13").Value = 0.287
$ws.Range("BO13").Value = 0.536
$ws.Range("BP13").Value = 0.795
$ws.Range("BQ13").Value = 0.734
